# Update NATMI LR-pair metrics per Dr Hou's advice (Ligand-expressing cells / Receptor-expressing cells corrected from 1 to 3)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map column letters to column indices
$colMap = @{ "E" = 5; "G" = 7; "H" = 8; "I" = 9; "J" = 10; "K" = 11; "M" = 13; "N" = 14; "O" = 15; "P" = 16; "Q" = 17; "R" = 18; "S" = 19; "T" = 20 }

$rowsData = @(
    @{ Row = 2; Values = @{ "E" = 3; "G" = 130.955829; "H" = 392.867487; "I" = 0.5336535908353144; "J" = 0.5336535908353144; "K" = 3; "M" = 224.2321046666667; "N" = 672.696314; "O" = 0.9009864013525987; "P" = 0.9009864013525988; "Q" = 29364.5011550381; "R" = 264280.5103953429; "S" = 0.480814628375602; "T" = 0.4808146283756021 } },
    @{ Row = 3; Values = @{ "E" = 3; "G" = 130.955829; "H" = 392.867487; "I" = 0.5336535908353144; "J" = 0.5336535908353144; "K" = 3; "M" = 1.023704333333333; "N" = 3.071113; "O" = 0.004113343558497904; "P" = 0.004113343558497904; "Q" = 134.060049622559; "R" = 1206.540446603031; "S" = 0.002195100560331716; "T" = 0.002195100560331717 } },
    @{ Row = 4; Values = @{ "E" = 3; "G" = 130.955829; "H" = 392.867487; "I" = 0.5336535908353144; "J" = 0.5336535908353144; "K" = 3; "M" = 23.61820766666667; "N" = 70.854623; "O" = 0.0949002550889034; "P" = 0.09490025508890343; "Q" = 3092.941964482489; "R" = 27836.4776803424; "S" = 0.05064386189938062; "T" = 0.05064386189938063 } },
    @{ Row = 5; Values = @{ "E" = 3; "G" = 66.39541; "H" = 199.18623; "I" = 0.2705656497465488; "J" = 0.2705656497465488; "K" = 3; "M" = 224.2321046666667; "N" = 672.696314; "O" = 0.9009864013525987; "P" = 0.9009864013525988; "Q" = 14887.98252450625; "R" = 133991.8427205562; "S" = 0.2437759710947707; "T" = 0.2437759710947707 } },
    @{ Row = 6; Values = @{ "E" = 3; "G" = 66.39541; "H" = 199.18623; "I" = 0.2705656497465488; "J" = 0.2705656497465488; "K" = 3; "M" = 1.023704333333333; "N" = 3.071113; "O" = 0.004113343558497904; "P" = 0.004113343558497904; "Q" = 67.96926893044333; "R" = 611.7234203739899; "S" = 0.001112929472535767; "T" = 0.001112929472535767 } },
    @{ Row = 7; Values = @{ "E" = 3; "G" = 66.39541; "H" = 199.18623; "I" = 0.2705656497465488; "J" = 0.2705656497465488; "K" = 3; "M" = 23.61820766666667; "N" = 70.854623; "O" = 0.0949002550889034; "P" = 0.09490025508890343; "Q" = 1568.140581493477; "R" = 14113.26523344129; "S" = 0.02567674917924238; "T" = 0.02567674917924238 } },
    @{ Row = 8; Values = @{ "E" = 3; "G" = 48.043585; "H" = 144.130755; "I" = 0.1957807594181367; "J" = 0.1957807594181367; "K" = 3; "M" = 224.2321046666667; "N" = 672.696314; "O" = 0.9009864013525987; "P" = 0.9009864013525988; "Q" = 10772.9141802819; "R" = 96956.22762253707; "S" = 0.1763958018822259; "T" = 0.1763958018822259 } },
    @{ Row = 9; Values = @{ "E" = 3; "G" = 48.043585; "H" = 144.130755; "I" = 0.1957807594181367; "J" = 0.1957807594181367; "K" = 3; "M" = 1.023704333333333; "N" = 3.071113; "O" = 0.004113343558497904; "P" = 0.004113343558497904; "Q" = 49.18242615336833; "R" = 442.641835380315; "S" = 0.0008053135256304205; "T" = 0.0008053135256304205 } },
    @{ Row = 10; Values = @{ "E" = 3; "G" = 48.043585; "H" = 144.130755; "I" = 0.1957807594181367; "J" = 0.1957807594181367; "K" = 3; "M" = 23.61820766666667; "N" = 70.854623; "O" = 0.0949002550889034; "P" = 0.09490025508890343; "Q" = 1134.703367581152; "R" = 10212.33030823036; "S" = 0.0185796440102804; "T" = 0.0185796440102804 } }
)

foreach ($rd in $rowsData) {
    $r = $rd.Row
    foreach ($col in $rd.Values.Keys) {
        $c = $colMap[$col]
        $ws.Cells.Item($r, $c).Value = $rd.Values[$col]
    }
}
